# Generate Report for Handback
# Refresh the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps for the first data
# row (row 2) on both the zh-cn and de-de report sheets, simulating a
# re-run of the handback status report generator.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-22 02:42:14"
$zhcn.Range("H2").Value = "2016-03-22 02:42:35"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-22 02:42:18"
$dede.Range("H2").Value = "2016-03-22 02:42:41"
